$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ---- ALC ----
$ws_ALC.Range("H98").Value = 26927910
$ws_ALC.Range("I98").Value = 8700876
$ws_ALC.Range("J98").Value = 166668500
$ws_ALC.Range("K98").Value = 8700876
$ws_ALC.Range("L98").Value = 166668500
$ws_ALC.Range("M98").Value = -8699378
$ws_ALC.Range("N98").Value = -166671496
$ws_ALC.Range("H122").Value = 26927910
$ws_ALC.Range("I122").Value = 8700876
$ws_ALC.Range("J122").Value = 166668500
$ws_ALC.Range("K122").Value = 26102628
$ws_ALC.Range("L122").Value = 500005500
$ws_ALC.Range("M122").Value = -26100178
$ws_ALC.Range("N122").Value = -500010400
$ws_ALC.Range("H134").Value = 226250
$ws_ALC.Range("J134").Value = 226250
$ws_ALC.Range("L134").Value = 226250
$ws_ALC.Range("N134").Value = -236390
$ws_ALC.Range("H135").Value = 62503876
$ws_ALC.Range("I135").Value = 2433.3333
$ws_ALC.Range("J135").Value = 100004744
$ws_ALC.Range("K135").Value = 21899.9997
$ws_ALC.Range("L135").Value = 900042696
$ws_ALC.Range("M135").Value = -19364.9997
$ws_ALC.Range("N135").Value = -900047766
$ws_ALC.Range("H136").Value = 36765
$ws_ALC.Range("J136").Value = 36765
$ws_ALC.Range("L136").Value = 36765
$ws_ALC.Range("N136").Value = -46965
$ws_ALC.Range("H139").Value = 90754.28999999999
$ws_ALC.Range("J139").Value = 90754.28999999999
$ws_ALC.Range("L139").Value = 90754.28999999999
$ws_ALC.Range("N139").Value = -101034.29
$ws_ALC.Range("H140").Value = 67300
$ws_ALC.Range("J140").Value = 67300
$ws_ALC.Range("L140").Value = 67300
$ws_ALC.Range("N140").Value = -77660

# ---- ARM ----
$ws_ARM.Range("H32").Value = 7362209
$ws_ARM.Range("I32").Value = 8124.3335
$ws_ARM.Range("J32").Value = 35727964
$ws_ARM.Range("K32").Value = 8124.3335
$ws_ARM.Range("L32").Value = 35727964
$ws_ARM.Range("M32").Value = -7837.3335
$ws_ARM.Range("N32").Value = -35728538
$ws_ARM.Range("H122").Value = 1265.4
$ws_ARM.Range("I122").Value = 807.55554
$ws_ARM.Range("K122").Value = 2422.66662
$ws_ARM.Range("M122").Value = 27.33338000000003

# ---- BSM ----
$ws_BSM.Range("H22").Value = 479.6842
$ws_BSM.Range("I22").Value = 500.77777
$ws_BSM.Range("J22").Value = 100
$ws_BSM.Range("K22").Value = 500.77777
$ws_BSM.Range("L22").Value = 100
$ws_BSM.Range("M22").Value = -327.77777
$ws_BSM.Range("N22").Value = -446
$ws_BSM.Range("H86").Value = 2697.0403
$ws_BSM.Range("I86").Value = 2000
$ws_BSM.Range("J86").Value = 36503.5
$ws_BSM.Range("K86").Value = 2000
$ws_BSM.Range("L86").Value = 36503.5
$ws_BSM.Range("M86").Value = -877
$ws_BSM.Range("N86").Value = -38749.5
$ws_BSM.Range("H89").Value = 2697.0403
$ws_BSM.Range("I89").Value = 2000
$ws_BSM.Range("J89").Value = 36503.5
$ws_BSM.Range("K89").Value = 10000
$ws_BSM.Range("L89").Value = 182517.5
$ws_BSM.Range("M89").Value = -4384
$ws_BSM.Range("N89").Value = -193749.5
$ws_BSM.Range("H126").Value = 0
$ws_BSM.Range("J126").Value = 0
$ws_BSM.Range("L126").Value = 0
$ws_BSM.Range("N126").Value = ""
$ws_BSM.Range("H128").Value = 1800
$ws_BSM.Range("I128").Value = 1800
$ws_BSM.Range("K128").Value = 5400
$ws_BSM.Range("M128").Value = -2910
$ws_BSM.Range("H130").Value = 3000000
$ws_BSM.Range("J130").Value = 3000000
$ws_BSM.Range("L130").Value = 3000000
$ws_BSM.Range("N130").Value = -3010040

# ---- CRP ----
$ws_CRP.Range("H39").Value = 5000
$ws_CRP.Range("I39").Value = 5000
$ws_CRP.Range("J39").Value = 0
$ws_CRP.Range("K39").Value = 5000
$ws_CRP.Range("L39").Value = 0
$ws_CRP.Range("M39").Value = -4609
$ws_CRP.Range("N39").Value = ""
$ws_CRP.Range("H41").Value = 6666.6665
$ws_CRP.Range("I41").Value = 5166.6665
$ws_CRP.Range("J41").Value = 8166.6665
$ws_CRP.Range("K41").Value = 5166.6665
$ws_CRP.Range("L41").Value = 8166.6665
$ws_CRP.Range("M41").Value = -4738.6665
$ws_CRP.Range("N41").Value = -9022.666499999999
$ws_CRP.Range("H43").Value = 18828.572
$ws_CRP.Range("J43").Value = 18828.572
$ws_CRP.Range("L43").Value = 18828.572
$ws_CRP.Range("N43").Value = -19196.572
$ws_CRP.Range("H49").Value = 5000
$ws_CRP.Range("I49").Value = 5000
$ws_CRP.Range("J49").Value = 0
$ws_CRP.Range("K49").Value = 5000
$ws_CRP.Range("L49").Value = 0
$ws_CRP.Range("M49").Value = -4818
$ws_CRP.Range("N49").Value = ""
$ws_CRP.Range("H50").Value = 13575.571
$ws_CRP.Range("J50").Value = 13575.571
$ws_CRP.Range("L50").Value = 13575.571
$ws_CRP.Range("N50").Value = -14825.571
$ws_CRP.Range("H51").Value = 9400.4
$ws_CRP.Range("J51").Value = 9400.4
$ws_CRP.Range("L51").Value = 9400.4
$ws_CRP.Range("N51").Value = -10872.4
$ws_CRP.Range("H59").Value = 15303.454
$ws_CRP.Range("J59").Value = 16055.8
$ws_CRP.Range("L59").Value = 16055.8
$ws_CRP.Range("N59").Value = -18345.8
$ws_CRP.Range("H60").Value = 7562.8335
$ws_CRP.Range("J60").Value = 8099.4
$ws_CRP.Range("L60").Value = 8099.4
$ws_CRP.Range("N60").Value = -9121.4
$ws_CRP.Range("H61").Value = 9400.4
$ws_CRP.Range("J61").Value = 9400.4
$ws_CRP.Range("L61").Value = 9400.4
$ws_CRP.Range("N61").Value = -10096.4
$ws_CRP.Range("H74").Value = 17409.572
$ws_CRP.Range("J74").Value = 18677.834
$ws_CRP.Range("L74").Value = 18677.834
$ws_CRP.Range("N74").Value = -20425.834
$ws_CRP.Range("H76").Value = 20000
$ws_CRP.Range("I76").Value = 20000
$ws_CRP.Range("K76").Value = 20000
$ws_CRP.Range("M76").Value = -19685
$ws_CRP.Range("H77").Value = 17409.572
$ws_CRP.Range("J77").Value = 18677.834
$ws_CRP.Range("L77").Value = 56033.50199999999
$ws_CRP.Range("N77").Value = -64769.50199999999
$ws_CRP.Range("H79").Value = 20000
$ws_CRP.Range("I79").Value = 20000
$ws_CRP.Range("K79").Value = 20000
$ws_CRP.Range("M79").Value = -18908
$ws_CRP.Range("H101").Value = 18828.572
$ws_CRP.Range("J101").Value = 18828.572
$ws_CRP.Range("L101").Value = 18828.572
$ws_CRP.Range("N101").Value = -25318.572
$ws_CRP.Range("H122").Value = 5960.2104
$ws_CRP.Range("I122").Value = 7499.5713
$ws_CRP.Range("J122").Value = 1650
$ws_CRP.Range("K122").Value = 22498.7139
$ws_CRP.Range("L122").Value = 4950
$ws_CRP.Range("M122").Value = -20048.7139
$ws_CRP.Range("N122").Value = -9850

# ---- CUL ----
$ws_CUL.Range("H5").Value = 3984210.8
$ws_CUL.Range("I5").Value = 3205632.2
$ws_CUL.Range("J5").Value = 5682927.5
$ws_CUL.Range("K5").Value = 9616896.600000001
$ws_CUL.Range("L5").Value = 17048782.5
$ws_CUL.Range("M5").Value = -9616784.600000001
$ws_CUL.Range("N5").Value = -17049006.5
$ws_CUL.Range("H122").Value = 1076.238
$ws_CUL.Range("I122").Value = 290.7
$ws_CUL.Range("J122").Value = 1790.3636
$ws_CUL.Range("K122").Value = 2616.3
$ws_CUL.Range("L122").Value = 16113.2724
$ws_CUL.Range("M122").Value = -166.2999999999997
$ws_CUL.Range("N122").Value = -21013.2724
$ws_CUL.Range("H135").Value = 3984210.8
$ws_CUL.Range("I135").Value = 3205632.2
$ws_CUL.Range("J135").Value = 5682927.5
$ws_CUL.Range("K135").Value = 28850689.8
$ws_CUL.Range("L135").Value = 51146347.5
$ws_CUL.Range("M135").Value = -28848154.8
$ws_CUL.Range("N135").Value = -51151417.5

# ---- GSM ----
$ws_GSM.Range("H102").Value = 2941.415
$ws_GSM.Range("I102").Value = 3115.8164
$ws_GSM.Range("J102").Value = 805
$ws_GSM.Range("K102").Value = 3115.8164
$ws_GSM.Range("L102").Value = 805
$ws_GSM.Range("M102").Value = -1493.8164
$ws_GSM.Range("N102").Value = -4049
$ws_GSM.Range("H122").Value = 6453873.5
$ws_GSM.Range("I122").Value = 45361.72
$ws_GSM.Range("J122").Value = 166666670
$ws_GSM.Range("K122").Value = 136085.16
$ws_GSM.Range("L122").Value = 500000010
$ws_GSM.Range("M122").Value = -133635.16
$ws_GSM.Range("N122").Value = -500004910

# ---- LTW ----
$ws_LTW.Range("H30").Value = 989
$ws_LTW.Range("I30").Value = 783.2
$ws_LTW.Range("J30").Value = 2018
$ws_LTW.Range("K30").Value = 783.2
$ws_LTW.Range("L30").Value = 2018
$ws_LTW.Range("M30").Value = -675.2
$ws_LTW.Range("N30").Value = -2234
$ws_LTW.Range("H35").Value = 5875
$ws_LTW.Range("I35").Value = 1166.6666
$ws_LTW.Range("J35").Value = 20000
$ws_LTW.Range("K35").Value = 1166.6666
$ws_LTW.Range("L35").Value = 20000
$ws_LTW.Range("M35").Value = -830.6666
$ws_LTW.Range("N35").Value = -20672
$ws_LTW.Range("H107").Value = 20000
$ws_LTW.Range("I107").Value = 20000
$ws_LTW.Range("K107").Value = 20000
$ws_LTW.Range("M107").Value = -18080

# ---- WVR ----
$ws_WVR.Range("H128").Value = 54617.5
$ws_WVR.Range("J128").Value = 54617.5
$ws_WVR.Range("L128").Value = 54617.5
$ws_WVR.Range("N128").Value = -64577.5
